# feat: add 2022-Q3 data
#
# 1. Insert a brand-new worksheet named "2022-Q3" right after the "总计"
#    sheet (i.e. before the existing "2022-Q1" sheet) and fill it with the
#    quarterly fund-holding table for 2022-Q3.
# 2. On the "总计" summary sheet, insert a new row for "2022-Q3" above the
#    existing "2022-Q1" row (which shifts down and has its running index
#    bumped from 0 to 1).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the new "2022-Q3" worksheet, positioned right after "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q3Sheet.Name = "2022-Q3"

# Copy the header formatting (bold font + thin border, centered) from the
# "总计" sheet's own header row so the new sheet's header reuses the same
# style, then overwrite the header labels.
$totalSheet.Range("B1:D1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)
$totalSheet.Range("A2").Copy()
$q3Sheet.Range("A2").PasteSpecial(-4122)

$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

$q3Sheet.Range("A2").Value = 0

$q3Sheet.Range("B2").NumberFormat = "@"
$q3Sheet.Range("B2").Value = "014294"
$q3Sheet.Range("B2").ClearFormats()

$q3Sheet.Range("C2").NumberFormat = "@"
$q3Sheet.Range("C2").Value = "南方北交所精选两年定开混合"
$q3Sheet.Range("C2").ClearFormats()

$q3Sheet.Range("D2").NumberFormat = "@"
$q3Sheet.Range("D2").Value = "4.26"
$q3Sheet.Range("D2").ClearFormats()

$q3Sheet.Range("E2").NumberFormat = "@"
$q3Sheet.Range("E2").Value = "75.23"
$q3Sheet.Range("E2").ClearFormats()

$q3Sheet.Range("F2").NumberFormat = "@"
$q3Sheet.Range("F2").Value = "1.66"
$q3Sheet.Range("F2").ClearFormats()

$q3Sheet.Range("G2").NumberFormat = "@"
$q3Sheet.Range("G2").Value = "0.0707"
$q3Sheet.Range("G2").ClearFormats()

$q3Sheet.Range("H2").Value = 9

# ---------------------------------------------------------------------
# Step 2: insert the "2022-Q3" summary row on "总计", above "2022-Q1"
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# Give the new A2 the same style as the (shifted) A3 cell (bold/border/
# centered numeric style used throughout this column).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.07000000000000001

# The pre-existing "2022-Q1" row (now row 3) keeps its values, but its
# running index bumps from 0 to 1.
$totalSheet.Range("A3").Value = 1
